$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$newSystemText = "SYSTEM Exibe os detalhes relativos àquela prestação de contas (nome do beneficiário, dados básicos da solicitação e documentos anexos); Exibe o histórico da tramitação da prestação de contas."
$newChefeText = "Chefe Verifica o histório da tramitação da prestação de contas e clica para analisar a prestação de contas."

# Cells holding the "SYSTEM Exibe os detalhes..." text (shared string reused across the 4 test cases)
$systemCells = @("D11", "D23", "D34", "D53")
foreach ($addr in $systemCells) {
    $ws.Range($addr).Value = $newSystemText
}

# Cells holding the "Chefe Clica para analisar..." text (shared string reused across the 4 test cases)
$chefeCells = @("B12", "B24", "B35", "B54")
foreach ($addr in $chefeCells) {
    $ws.Range($addr).Value = $newChefeText
}
